$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H43").Value = 2991.5217
$ws.Range("I43").Value = 651.1111
$ws.Range("J43").Value = 4496.0713
$ws.Range("K43").Value = 651.1111
$ws.Range("L43").Value = 4496.0713
$ws.Range("M43").Value = -582.1111
$ws.Range("N43").Value = -4634.0713
$ws.Range("H80").Value = 12182917
$ws.Range("I80").Value = 462.75
$ws.Range("J80").Value = 20304554
$ws.Range("K80").Value = 1388.25
$ws.Range("L80").Value = 60913662
$ws.Range("M80").Value = -390.25
$ws.Range("N80").Value = -60915658
$ws.Range("H83").Value = 12182917
$ws.Range("I83").Value = 462.75
$ws.Range("J83").Value = 20304554
$ws.Range("K83").Value = 4164.75
$ws.Range("L83").Value = 182740986
$ws.Range("M83").Value = 827.25
$ws.Range("N83").Value = -182750970
$ws.Range("H92").Value = 946.75
$ws.Range("I92").Value = 946.75
$ws.Range("K92").Value = 946.75
$ws.Range("M92").Value = 301.25
$ws.Range("H123").Value = 24990
$ws.Range("J123").Value = 24990
$ws.Range("L123").Value = 24990
$ws.Range("N123").Value = -34790
$ws.Range("H129").Value = 1557.5883
$ws.Range("I129").Value = 135.66667
$ws.Range("J129").Value = 1695.1936
$ws.Range("K129").Value = 407.00001
$ws.Range("L129").Value = 5085.5808
$ws.Range("M129").Value = 4592.99999
$ws.Range("N129").Value = -15085.5808
$ws.Range("H132").Value = 47623330
$ws.Range("I132").Value = 50004420
$ws.Range("J132").Value = 1500
$ws.Range("K132").Value = 150013260
$ws.Range("L132").Value = 4500
$ws.Range("M132").Value = -150010730
$ws.Range("N132").Value = -9560
$ws.Range("H134").Value = 0
$ws.Range("J134").Value = 0
$ws.Range("L134").Value = 0
$ws.Range("N134").ClearContents()
$ws.Range("H135").Value = 15155133
$ws.Range("I135").Value = 634.7931
$ws.Range("K135").Value = 5713.1379
$ws.Range("M135").Value = -3178.1379
$ws.Range("H137").Value = 94694.72
$ws.Range("I137").Value = 115761.83
$ws.Range("J137").Value = 2526.125
$ws.Range("K137").Value = 347285.49
$ws.Range("L137").Value = 7578.375
$ws.Range("M137").Value = -344735.49
$ws.Range("N137").Value = -12678.375

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 7938618.5
$ws.Range("I61").Value = 10418138
$ws.Range("J61").Value = 4155.3
$ws.Range("K61").Value = 10418138
$ws.Range("L61").Value = 4155.3
$ws.Range("M61").Value = -10417926
$ws.Range("N61").Value = -4579.3
$ws.Range("H136").Value = 7938618.5
$ws.Range("I136").Value = 10418138
$ws.Range("J136").Value = 4155.3
$ws.Range("K136").Value = 31254414
$ws.Range("L136").Value = 12465.9
$ws.Range("M136").Value = -31251864
$ws.Range("N136").Value = -17565.9

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 959.7586
$ws.Range("I94").Value = 557.58826
$ws.Range("J94").Value = 1529.5
$ws.Range("K94").Value = 557.58826
$ws.Range("L94").Value = 1529.5
$ws.Range("M94").Value = -106.58826
$ws.Range("N94").Value = -2431.5
$ws.Range("H134").Value = 3760.15
$ws.Range("I134").Value = 3827.1892
$ws.Range("J134").Value = 2933.3333
$ws.Range("K134").Value = 11481.5676
$ws.Range("L134").Value = 8799.999899999999
$ws.Range("M134").Value = -8946.567599999998
$ws.Range("N134").Value = -13869.9999

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H52").Value = 20661
$ws.Range("I52").Value = 7709
$ws.Range("J52").Value = 22280
$ws.Range("K52").Value = 7709
$ws.Range("L52").Value = 22280
$ws.Range("M52").Value = -7415
$ws.Range("N52").Value = -22868

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 1219.5416
$ws.Range("J5").Value = 1839.3158
$ws.Range("L5").Value = 5517.9474
$ws.Range("N5").Value = -5741.9474
$ws.Range("H33").Value = 45.3
$ws.Range("I33").Value = 15
$ws.Range("J33").Value = 65.5
$ws.Range("K33").Value = 90
$ws.Range("L33").Value = 393
$ws.Range("M33").Value = 193
$ws.Range("N33").Value = -959
$ws.Range("H80").Value = 22900.2
$ws.Range("J80").Value = 27875.25
$ws.Range("L80").Value = 83625.75
$ws.Range("N80").Value = -85497.75
$ws.Range("H83").Value = 22900.2
$ws.Range("J83").Value = 27875.25
$ws.Range("L83").Value = 250877.25
$ws.Range("N83").Value = -260237.25
$ws.Range("H129").Value = 191950.9
$ws.Range("I129").Value = 503.6
$ws.Range("J129").Value = 365993.9
$ws.Range("K129").Value = 1510.8
$ws.Range("L129").Value = 1097981.7
$ws.Range("M129").Value = 3489.2
$ws.Range("N129").Value = -1107981.7
$ws.Range("H131").Value = 681.4433
$ws.Range("J131").Value = 721.3563
$ws.Range("L131").Value = 2164.0689
$ws.Range("N131").Value = -12244.0689
$ws.Range("H135").Value = 1219.5416
$ws.Range("J135").Value = 1839.3158
$ws.Range("L135").Value = 16553.8422
$ws.Range("N135").Value = -21623.8422
$ws.Range("H137").Value = 13893467
$ws.Range("I137").Value = 1294.2858
$ws.Range("J137").Value = 19613772
$ws.Range("K137").Value = 3882.8574
$ws.Range("L137").Value = 58841316
$ws.Range("M137").Value = 1217.1426
$ws.Range("N137").Value = -58851516

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H46").Value = 34866.668
$ws.Range("J46").Value = 34866.668
$ws.Range("L46").Value = 34866.668
$ws.Range("N46").Value = -35178.668
$ws.Range("H97").Value = 6005.5
$ws.Range("I97").Value = 0
$ws.Range("J97").Value = 6005.5
$ws.Range("K97").Value = 0
$ws.Range("L97").Value = 6005.5
$ws.Range("M97").ClearContents()
$ws.Range("N97").Value = -6997.5
$ws.Range("H126").Value = 3466.1482
$ws.Range("I126").Value = 2442.4707
$ws.Range("K126").Value = 7327.4121
$ws.Range("M126").Value = -4857.4121
$ws.Range("H132").Value = 4254308
$ws.Range("I132").Value = 6688870
$ws.Range("K132").Value = 20066610
$ws.Range("M132").Value = -20064080

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4532.3335
$ws.Range("I7").Value = 3898.8
$ws.Range("K7").Value = 3898.8
$ws.Range("M7").Value = -3786.8
$ws.Range("H16").Value = 804.0714
$ws.Range("I16").Value = 805.8889
$ws.Range("K16").Value = 805.8889
$ws.Range("M16").Value = -635.8889
$ws.Range("H40").Value = 4620.5264
$ws.Range("I40").Value = 4377.222
$ws.Range("J40").Value = 9000
$ws.Range("K40").Value = 4377.222
$ws.Range("L40").Value = 9000
$ws.Range("M40").Value = -4241.222
$ws.Range("N40").Value = -9272
$ws.Range("H126").Value = 4532.3335
$ws.Range("I126").Value = 3898.8
$ws.Range("K126").Value = 11696.4
$ws.Range("M126").Value = -9226.400000000001
$ws.Range("H132").Value = 276776.53
$ws.Range("I132").Value = 448839.72
$ws.Range("J132").Value = 3499.7058
$ws.Range("K132").Value = 1346519.16
$ws.Range("L132").Value = 10499.1174
$ws.Range("M132").Value = -1343989.16
$ws.Range("N132").Value = -15559.1174

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 1912
$ws.Range("J96").Value = 2596
$ws.Range("L96").Value = 2596
$ws.Range("N96").Value = -5342
$ws.Range("H126").Value = 2429.6843
$ws.Range("J126").Value = 3066.2856
$ws.Range("L126").Value = 9198.856800000001
$ws.Range("N126").Value = -14138.8568
$ws.Range("H132").Value = 12195994
$ws.Range("I132").Value = 14286313
$ws.Range("K132").Value = 42858939
$ws.Range("M132").Value = -42856409

Write-Output "done"